# "Fix up onto github triggers"
#
# The "components" sheet's column D held the repo-root GitHub URLs that are
# used as webhook/trigger URLs; they had a trailing slash which broke the
# trigger matching. Strip the trailing slash from all seven of them and make
# sure the cells are real hyperlinks (reusing the same "Hyperlink" style
# already used elsewhere on the sheet). Also make "components" the active
# sheet (it was "clusters").

$wb = $excel.ActiveWorkbook

$clusters = $wb.Worksheets.Item("clusters")
$components = $wb.Worksheets.Item("components")
$environments = $wb.Worksheets.Item("environments")

# New (slash-stripped) target URLs for column D, rows 1..7.
$newUrls = @{
    1 = "https://github.com/ConductionNL/environment-component"
    2 = "https://github.com/ConductionNL/user-component"
    3 = "https://github.com/ConductionNL/Commongroundregistratiecomponent"
    4 = "https://github.com/ConductionNL/procesregistratiecomponent"
    5 = "https://github.com/ConductionNL/medewerkercatalogus"
    6 = "https://github.com/ConductionNL/webresourcecatalogus"
    7 = "https://github.com/ConductionNL/digispoof-interface"
}

foreach ($row in 1..7) {
    $components.Cells.Item($row, 4).Value = $newUrls[$row]
}

# Hyperlink.Delete()/per-item deletion is a no-op in this host, and a
# single-cell Range.Hyperlinks.Delete() wipes the whole sheet's collection
# instead of just that cell's link - so rebuild the full set from a clean
# slate, re-adding the two that already existed (E4, D7) first (in their
# original order) and then the six brand-new ones, which reproduces the
# same rId1..rId8 relationship order the links would naturally get.
$components.Cells.Hyperlinks.Delete()

$components.Hyperlinks.Add($components.Range("E4"), "https://github.com/ConductionNL/procesregistratiecomponent/api/helm") | Out-Null
$components.Hyperlinks.Add($components.Range("D7"), $newUrls[7]) | Out-Null
$components.Hyperlinks.Add($components.Range("D1"), $newUrls[1]) | Out-Null
$components.Hyperlinks.Add($components.Range("D2"), $newUrls[2]) | Out-Null
$components.Hyperlinks.Add($components.Range("D3"), $newUrls[3]) | Out-Null
$components.Hyperlinks.Add($components.Range("D4"), $newUrls[4]) | Out-Null
$components.Hyperlinks.Add($components.Range("D5"), $newUrls[5]) | Out-Null
$components.Hyperlinks.Add($components.Range("D6"), $newUrls[6]) | Out-Null

# Hyperlinks.Add() drops a brand-new duplicate style on every cell it
# touches; re-apply the real shared "Hyperlink" cell style (underline +
# theme color, same one already used on E4/D7) so every linked cell in
# column D looks consistent again.
$components.Range("D1:D7").Style = "Hyperlink"
$components.Range("E4").Style = "Hyperlink"

# components becomes the active/selected sheet (clusters loses
# tabSelected); scroll/select near the right edge of the used range.
$components.Activate()
$components.Range("G1").Select() | Out-Null
